# ----------------------------------------------------------------------------
# Applies three related edits to "Instructions on how to run HH models.docx":
#
#  1. "507" -> "530" in "Each run will be stopped after 507 secs ("507000")."
#     The resulting sentence ends up split across three runs (same formatting):
#       " domain ("5") 31 times ("31"). Ea" | "ch run will be stopped after 530
#       secs ("530" | "000")."
#
#  2. The "_GoBack" bookmark (left behind after the previous edit, sitting
#     right after "...stops after 346") is removed from there.
#
#  3. The "_GoBack" bookmark re-appears between "Th" and "anks" inside
#     "Thanks for showing interest in our research", splitting that run in
#     two (same formatting on both halves).
#
# Note: in this engine, assigning Range.Text anywhere inside a paragraph
# re-normalises (merges) every same-formatting run in that paragraph. So we
# do the text replacement first (accepting the temporary merge) and then
# recreate every needed run boundary - including ones that already existed,
# like the boundary around "VRP" - via a harmless Bold toggle (set then
# immediately clear), which *does* preserve a hard split at that offset.
# ----------------------------------------------------------------------------

$d = $word.ActiveDocument

# ========== Change 1: "507" -> "530" ==========

$full = $d.Content.Text
$sentence = 'Each run will be stopped after 507 secs ("507000").'
$idxEach = $full.IndexOf($sentence)
if ($idxEach -lt 0) { throw "Could not find target sentence" }

$idx507_1 = $full.IndexOf('507', $idxEach)
$idx507_2 = $full.IndexOf('507', $idx507_1 + 3)

$d.Range($idx507_1, $idx507_1 + 3).Text = "530"
$d.Range($idx507_2, $idx507_2 + 3).Text = "530"

# Recompute offsets in the now-merged paragraph text.
$full2 = $d.Content.Text
$newSentence = 'Each run will be stopped after 530 secs ("530000").'
$idxEach2 = $full2.IndexOf($newSentence)
$endOfSentence = $idxEach2 + $newSentence.Length

$idxVRP = $full2.IndexOf('VRP')
$splitVRPstart = $idxVRP            # boundary before "VRP"  (after "...of the ")
$splitVRPend   = $idxVRP + 3        # boundary after "VRP"   (before " domain")

$splitA = $idxEach2 + 2             # boundary between "Ea" and "ch"

$idx530_2 = $full2.IndexOf('530', $idxEach2 + 35)
$splitB = $idx530_2 + 3             # boundary between "530" and "000"

# Recreate boundaries from rightmost to leftmost so earlier offsets stay valid.
foreach ($pt in @($splitB, $splitA, $splitVRPend, $splitVRPstart)) {
    $rMark = $d.Range($pt, $endOfSentence)
    $rMark.Bold = 1
    $rMark.Bold = 0
}

# ========== Changes 2 & 3: move the "_GoBack" bookmark ==========
# Moving/re-adding a bookmark with the same name repositions it (and drops
# the old location), and also gives us the run split at "Th" | "anks" as a
# side effect, matching the target structure.

$full3 = $d.Content.Text
$idxThanks = $full3.IndexOf('Thanks for showing interest')
if ($idxThanks -lt 0) { throw "Could not find 'Thanks for showing interest'" }
$splitPoint = $idxThanks + 2   # between "Th" and "anks"

$d.Bookmarks.Add("_GoBack", $d.Range($splitPoint, $splitPoint))
